# Report 1 small update
# Adds a new paragraph right before the closing "Λήξη αναφοράς." paragraph,
# after the (empty) paragraph that follows the table.

$d = $word.ActiveDocument

# Locate the closing paragraph by its unique text. Find.Execute mutates the
# Range it is called on in place, so keep a handle to that same Range.
$rng = $d.Content
$found = $rng.Find.Execute("Λήξη αναφοράς.", $false, $false, $false, `
    $false, $false, $true, 1, $false, "", 0)

$closingPara = $rng.Paragraphs(1)

# The paragraph immediately before the closing one is the empty paragraph
# sitting right after the table - that's where the new text belongs.
$precedingEmptyPara = $closingPara.Previous()

# Insert a new paragraph after it; it inherits the (plain, non-bold,
# non-centered) formatting of the preceding empty paragraph.
$precedingEmptyPara.Range.InsertParagraphAfter()

# Grab the newly created paragraph and fill in its text.
$newPara = $precedingEmptyPara.Next()
$newPara.Range.Text = "Η δυσλειτουργία αυτή αν και κυρίως άκακη και σπάνια με " + `
    "ένα ποσοστό εμφάνισης μικρότερο του 1%, είναι αξιοσημείωτη καθώς " + `
    "μπορεί να μπερδέψει τον μελλοντικό χρήστη και να οδηγήσει σε " + `
    "περεταίρω προβλήματα."
